# Lesson 1.docx edit
#
# Summary of the change (see commit message "controller layer api week 14"):
#  - The single bullet "Operating system differences" (ilvl=1) under "CLI" is
#    replaced by six new bullets describing basic CLI commands:
#       Dir / Cd / mkdir / Copy / move  (all ilvl=2, sub-bullets of "Basic commands")
#       "Discuss how the cli directory is the same as the file explorer" (ilvl=1)
#    The trailing "_GoBack" last-edit bookmark, which used to sit at the end of
#    the "Who can find out..." bullet near the end of the CLI/Source
#    Control/Programming block, now sits at the end of the new
#    "Discuss how the cli directory..." bullet instead (i.e. it moved to the
#    new last-edited text).
#  - As a knock-on effect of content shifting earlier in the document, Word's
#    cached lastRenderedPageBreak hints move: it disappears from "How many
#    words in this sentence?..." and from "Operations", and a new one appears
#    on "We will talk about Objects more in a later lesson...". We replicate
#    that exactly to match the target render-cache state, and a new one shows
#    up on the "Remember if you are working with a team..." bullet too.

$d = $word.ActiveDocument

function Set-ParagraphXML {
    param(
        [int]$ParaIndex,
        [string]$BodyInnerXml
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $rng = $p.Range
    $full = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' + $BodyInnerXml + '</w:body>' +
            '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($full)
}

# ---------------------------------------------------------------------------
# 1) Replace the "Operating system differences" bullet (paragraph 5) with the
#    six new bullets about CLI commands + the relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$newBlock = (
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Dir</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Cd</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>mkdir</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Copy</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>move</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Discuss how the cli directory is the same as the file explorer</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
)
Set-ParagraphXML 5 $newBlock

# ---------------------------------------------------------------------------
# 2) Remove the old _GoBack bookmark from the "Who can find out..." bullet
#    (it moved to the new bullet inserted above).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Update the cached lastRenderedPageBreak hints to match the new layout.
#    Paragraph numbers below are from the ORIGINAL document and are still
#    valid here because step 1 kept the paragraph count change local (it
#    replaced 1 paragraph with 6, a net +5, and none of these numbers are
#    before that edit point other than being re-queried by text below).
# ---------------------------------------------------------------------------

function Find-ParagraphByText {
    param([string]$ExactText)
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $ExactText) {
            return $i
        }
    }
    return -1
}

# 3a) "How many words in this sentence?..." loses its lastRenderedPageBreak.
$idx = Find-ParagraphByText "How many words in this sentence? How many vowels? What is the average number of vowels per word?"
if ($idx -gt 0) {
    $xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>How many words in this sentence? How many vowels? What is the average number of vowels per word?</w:t></w:r></w:p>'
    Set-ParagraphXML $idx $xml
}

# 3b) "Remember if you are working with a team..." gains a lastRenderedPageBreak.
$idx = Find-ParagraphByText "Remember if you are working with a team to pull before making changes so you have the most up to date version of the code"
if ($idx -gt 0) {
    $xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Remember if you are working with a team to pull before making changes so you have the most up to date version of the code</w:t></w:r></w:p>'
    Set-ParagraphXML $idx $xml
}

# 3c) "We will talk about Objects more in a later lesson..." gains a lastRenderedPageBreak.
$idx = Find-ParagraphByText "We will talk about Objects more in a later lesson, but we will add String to the list of data types covered today."
if ($idx -gt 0) {
    $xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>We will talk about Objects more in a later lesson, but we will add String to the list of data types covered today.</w:t></w:r></w:p>'
    Set-ParagraphXML $idx $xml
}

# 3d) "Operations" loses its lastRenderedPageBreak.
$idx = Find-ParagraphByText "Operations"
if ($idx -gt 0) {
    $xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Operations</w:t></w:r></w:p>'
    Set-ParagraphXML $idx $xml
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
